# Apply cryptos list update (Tue Sep 26 23:35:18 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds numeric-looking price strings (e.g. "212.12") that must
# stay plain text, matching the source feed formatting - force text before assigning.
$ws.Range('D2').Value = '26.183.01'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.12'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.815.03'
$ws.Range('D13').Value = '1.599.45'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.64'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '26.183.05'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '213.79'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.33'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.67'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.09'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0494'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.17'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Value = '1.420.99'
$ws.Range('E33').Value = '  +7.74%  '
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.587'
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.90'
$ws.Range('E39').Value = '  +4.77%  '
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.988'
$ws.Range('E42').Value = '  -8.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.764'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = '1.726.57'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '60.95'
$ws.Range('E46').Value = '  -2.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '86.99'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0501'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0955'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('E51').Value = '  -0.29%  '
